# Update the public EPEX spot / Gaz / CO2 Excel workbook with the latest day
# of data (12-jul for "Prix Spot", 2025-07-10 for "Gaz" and "CO2").

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": append a new date column (AC) -----------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("AB1").Copy() | Out-Null
$wsSpot.Range("AC1").PasteSpecial(-4122) | Out-Null
$wsSpot.Range("AC1").Value = "12-jul"

$spotValues = @(100, 90.03, 82.67, 70.05, 63.95, 71.87, 69.45, 70.89, 65.6, 59.78, 18.19, 3, 0.22, 0, 0, 0.22, 13.68, 35, 50.08, 75.4, 107.75, 99.72, 117.17, 105.98)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 29).Value = $spotValues[$i]
}

# --- Sheet "Gaz": append the new day's row (26) ------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force the date-like string to be stored as plain text (not auto-converted to
# a date serial) while keeping the same (unstyled) look as the rows above it.
$wsGaz.Range("A26").NumberFormat = "@"
$wsGaz.Range("A26").Value = "2025-07-10"
$wsGaz.Range("A25").Copy() | Out-Null
$wsGaz.Range("A26").PasteSpecial(-4122) | Out-Null
$wsGaz.Range("B26").Value = 34.4

# --- Sheet "CO2": append the new day's row (26) ------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A26").NumberFormat = "@"
$wsCo2.Range("A26").Value = "2025-07-10"
$wsCo2.Range("A25").Copy() | Out-Null
$wsCo2.Range("A26").PasteSpecial(-4122) | Out-Null
$wsCo2.Range("B26").Value = 69.8
